$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Range("D2").Value = '43.450.71'
$ws.Range("E2").Value = '  +2.61%  '

# Row 3
$ws.Range("D3").Value = '2.405.94'
$ws.Range("E3").Value = '  +8.33%  '

# Row 4
$ws.Range("E4").Value = '  -0.19%  '

# Row 5
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = '325.16'
$ws.Range("E5").Value = '  +11.95%  '

# Row 6
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = '105.28'
$ws.Range("E6").Value = '  -5.38%  '

# Row 7
$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = '0.658'
$ws.Range("E7").Value = '  +5.32%  '

# Row 9
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = '0.661'
$ws.Range("E9").Value = '  +10.75%  '

# Row 10
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = '42.35'
$ws.Range("E10").Value = '  -3.05%  '

# Row 11
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = '0.0947'
$ws.Range("E11").Value = '  +4.24%  '

# Row 12
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = '8.62'
$ws.Range("E12").Value = '  +0.24%  '

# Row 13
$ws.Range("B13").Value = 'Polygon'
$ws.Range("C13").Value = 'https://coinranking.com/coin/uW2tk-ILY0ii+polygon-matic'
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = '1.02'
$ws.Range("E13").Value = '  +1.50%  '

# Row 14
$ws.Range("B14").Value = 'Chainlink'
$ws.Range("C14").Value = 'https://coinranking.com/coin/VLqpJwogdhHNb+chainlink-link'
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = '17.28'
$ws.Range("E14").Value = '  +16.19%  '

# Row 15
$ws.Range("E15").Value = '  +3.19%  '

# Row 16
$ws.Range("D16").Value = '2.770.07'
$ws.Range("E16").Value = '  +8.34%  '

# Row 17
$ws.Range("D17").Value = '2.406.16'
$ws.Range("E17").Value = '  +8.14%  '

# Row 18
$ws.Range("D18").Value = '43.456.24'
$ws.Range("E18").Value = '  +2.73%  '

# Row 19
$ws.Range("B19").Value = 'ShibaInu'
$ws.Range("C19").Value = 'https://coinranking.com/coin/xz24e0BjL+shibainu-shib'
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = '0.0000110'
$ws.Range("E19").Value = '  +4.62%  '

# Row 20
$ws.Range("B20").Value = 'Uniswap'
$ws.Range("C20").Value = 'https://coinranking.com/coin/_H5FVG9iW+uniswap-uni'
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = '7.45'
$ws.Range("E20").Value = '  +4.89%  '

# Row 21
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = '76.09'
$ws.Range("E21").Value = '  +4.46%  '

# Row 22
$ws.Range("B22").Value = 'BitcoinCash'
$ws.Range("C22").Value = 'https://coinranking.com/coin/ZlZpzOJo43mIo+bitcoincash-bch'
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = '272.70'
$ws.Range("E22").Value = '  +15.36%  '

# Row 23
$ws.Range("B23").Value = 'PancakeSwap'
$ws.Range("C23").Value = 'https://coinranking.com/coin/ncYFcP709+pancakeswap-cake'
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = '3.48'
$ws.Range("E23").Value = '  +4.39%  '

# Row 24
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = '2.44'
$ws.Range("E24").Value = '  +1.88%  '

# Row 25
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = '9.74'
$ws.Range("E25").Value = '  +9.18%  '

# Row 26
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = '11.95'
$ws.Range("E26").Value = '  +5.01%  '

# Row 27
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = '1.00'
$ws.Range("E27").Value = '  +0.04%  '

# Row 28
$ws.Range("E28").Value = '  +0.04%  '

# Row 29
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = '22.98'
$ws.Range("E29").Value = '  +8.34%  '

# Row 30
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = '178.04'
$ws.Range("E30").Value = '  +2.95%  '

# Row 31
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = '2.21'
$ws.Range("E31").Value = '  +0.48%  '

# Row 32
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = '38.06'
$ws.Range("E32").Value = '  +1.24%  '

# Row 33
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = '3.23'
$ws.Range("E33").Value = '  +4.66%  '

# Row 34
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = '0.0937'
$ws.Range("E34").Value = '  +7.00%  '

# Row 35
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = '5.96'
$ws.Range("E35").Value = '  +6.59%  '

# Row 36
$ws.Range("E36").Value = '  +7.22%  '

# Row 37
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = '4.88'
$ws.Range("E37").Value = '  -1.75%  '

# Row 38
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = '4.10'
$ws.Range("E38").Value = '  -2.37%  '

# Row 39
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = '0.0371'
$ws.Range("E39").Value = '  -0.38%  '

# Row 40
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = '0.111'
$ws.Range("E40").Value = '  +6.91%  '

# Row 41
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = '2.86'
$ws.Range("E41").Value = '  +19.94%  '

# Row 42
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = '1.61'
$ws.Range("E42").Value = '  +22.86%  '

# Row 43
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = '127.46'
$ws.Range("E43").Value = '  +26.58%  '

# Row 44
$ws.Range("E44").Value = '  +1.94%  '

# Row 45
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = '70.01'
$ws.Range("E45").Value = '  -1.96%  '

# Row 46
$ws.Range("E46").Value = '  +4.04%  '

# Row 47
$ws.Range("E47").Value = '  +0.23%  '

# Row 48
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = '9.74'
$ws.Range("E48").Value = '  +15.96%  '

# Row 49
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = '5.69'
$ws.Range("E49").Value = '  +6.96%  '

# Row 50
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = '88.58'
$ws.Range("E50").Value = '  +59.45%  '

# Row 51
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = '1.32'
$ws.Range("E51").Value = '  +4.45%  '
